$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: add quantity 20 under "預計購買" (F6), and update unit-price / total (I6 / K6)
$ws.Range("F6").Value = 20
$ws.Range("I6").Value = "245(運費0)"
$ws.Range("K6").Value = "245*20=4900"

# Row 10: 杜邦簧片(母) -> 杜邦簧片(公), update price / total
$ws.Range("E10").Value = "杜邦簧片(公)"
$ws.Range("I10").Value = "3元/5個"
$ws.Range("K10").Value = "3*20=60"

# Row 11 (old 排針(公) row) is no longer needed - clear its contents entirely
$ws.Range("E11:K11").ClearContents()

# Row 12 (排針(母)): update quantity formula note, pack count, total
$ws.Range("F12").Value = "(19+19+4+8+4+2)*20"
$ws.Range("H12").Value = "40排"
$ws.Range("K12").Value = "4*40=160"

# New row 17: 杜邦塑膠殼(8pin)
$ws.Range("E17").Value = "杜邦塑膠殼(8pin)"
$ws.Range("F17").Value = 40
$ws.Range("I17").Value = "8元"
$ws.Range("H17").Value = "42個"
$ws.Range("K17").Value = "8*42=336"

# Row 16 (PCB): add new note in H16
$ws.Range("H16").Value = "25片"

# Remove the merges that used to span I6:J7 and K6:K7
$ws.Range("I6:J7").UnMerge()
$ws.Range("K6:K7").UnMerge()

# Update sheet view: scroll back to top-left and move selection
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("G3").Select()
